# Update Active_Outages.xlsx - 6/19/2025, 10:15:32 PM
# Refresh the "Elapsed Duration(Hrs)" figures across the per-region sheets
# and fix the stale Hub Site / Battery-Backup filler values on R1 row 6.

$wb = $excel.ActiveWorkbook

# --- R1 ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3959:29:33"
$ws.Range("G3").Value = "99:02:11"
$ws.Range("G4").Value = "122:02:11"
$ws.Range("D6").Value = "JED0925"
$ws.Range("J6").Value = "In progress"

# --- R2 ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12139:22:43"
$ws.Range("G3").Value = "3269:06:12"
$ws.Range("G4").Value = "507:17:46"

# --- R4 ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2985:12:32"
$ws.Range("G3").Value = "212:24:47"
$ws.Range("G4").Value = "100:37:12"
$ws.Range("G5").Value = "98:14:45"

# --- R5 ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "459:11:31"

# --- R6 ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "99:43:49"
